$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 2; $r -le $lastRow; $r++) {
    $typeVal = $ws.Cells.Item($r, 2).Value2
    $valVal = $ws.Cells.Item($r, 3).Value2

    if ([string]::IsNullOrEmpty($typeVal)) { continue }

    # Determine the sport category prefix (club / uil) from the old type value
    if ($typeVal -eq "club-sports") {
        $prefix = "club"
    } elseif ($typeVal -eq "uil-sports") {
        $prefix = "uil"
    } else {
        $prefix = $null
    }

    # Determine the gender suffix from the old value
    if ($valVal -eq "Baseball-Boys") {
        $suffix = "boys"
    } elseif ($valVal -eq "Baseball-Girls") {
        $suffix = "girls"
    } elseif ($valVal -eq "Baseball-Coed") {
        $suffix = "coed"
    } else {
        $suffix = $null
    }

    if ($prefix -ne $null -and $suffix -ne $null) {
        $ws.Cells.Item($r, 2).Value2 = "sports_" + $prefix + "_" + $suffix
        $ws.Cells.Item($r, 3).Value2 = "Baseball"
    }
}
